$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Through 2022-11-28" to "Through 2022-11-29"
$ws.Name = "Through 2022-11-29"

# Update the label in A12 to reflect new "through" date
$ws.Range("A12").Value = "November (through 11-29)"

# Update the November row (row 12) values
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = 76
$ws.Range("D12").Value = 105
$ws.Range("E12").Value = 67
$ws.Range("F12").Value = 48
$ws.Range("G12").Value = 205
$ws.Range("H12").Value = 192
$ws.Range("I12").Value = 114

# Update the Total row (row 13) values
$ws.Range("B13").Value = 290
$ws.Range("C13").Value = 562
$ws.Range("D13").Value = 815
$ws.Range("E13").Value = 682
$ws.Range("F13").Value = 530
$ws.Range("G13").Value = 1262
$ws.Range("H13").Value = 1633
$ws.Range("I13").Value = 1512
